# Apply the "edit Event / publish Event" permissions update:
#  - record the currently logged-in user ("admin") next to the header row
#  - append a new permission row (#11): HasPermissionToPublishEvent
#  - widen column B so the longer permission name fits
#  - leave the new row's cell selected, like a user who just typed it in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gains a value in column F ("admin"), with no special styling applied.
$ws.Range("F2").Value = "admin"

# Create row 10 by duplicating the formatting of the last existing data
# row (row 9), then filling in the new id/permission values.
$ws.Range("A9:B9").Copy($ws.Range("A10:B10"))
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "HasPermissionToPublishEvent"

# Widen column B to fit the newly added, longer permission name.
$ws.Columns.Item(2).ColumnWidth = 24.6

# Select the newly added cell, matching the saved cursor position.
$ws.Range("B10").Select() | Out-Null
